$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "BASE" to "Sheet1"
$ws.Name = "Sheet1"

# Data for rows 2-6 (row 1 is the header, already present)
$data = @(
    @(8000007, "Activo", 77,   "rrr",     "TS10", "MOTO", "SFE44R", 5465231, "A2", "2025/02/27", "2025/03/15", "2025/06/13"),
    @(8000008, "Activo", 88,   "tytt",    "TS11", "MOTO", "THV77G", 5465231, "A2", "2025/02/28", "2025/03/16", "2025/06/14"),
    @(8000010, "Activo", 1111, "saaa",    "TS13", "MOTO", "SFE44R", 5465231, "A2", "2025/03/02", "2025/03/18", "2025/06/16"),
    @(8000011, "Activo", 2222, "qwewqe",  "TS14", "MOTO", "THV77G", 5465231, "A2", "2025/03/03", "2025/03/19", "2025/06/17"),
    @(8000013, "Activo", 4444, "gjhgjhg", "TS16", "MOTO", "SFE44R", 5465231, "A2", "2025/03/05", "2025/03/21", "2025/06/19")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    # Mark the date-like text columns (J, K, L) as Text BEFORE writing so
    # Excel doesn't auto-convert the "2025/02/27" strings into date serials.
    $ws.Range($ws.Cells.Item($row, 10), $ws.Cells.Item($row, 12)).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
    $ws.Cells.Item($row, 12).Value = $vals[11]
    # Column M (OBSERVACIONES) is left blank, same as the source data - do
    # not touch it so any pre-existing empty-text state is preserved.
}
